# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" sheet (sheet1) and the "全部类型" sheet (sheet4), mirroring the
# data refresh captured in the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 97
$ws1.Range("F3").Value = 823
$ws1.Range("F6").Value = 136
$ws1.Range("F8").Value = 4856
$ws1.Range("F10").Value = 5176
$ws1.Range("F11").Value = 591
$ws1.Range("F12").Value = 1294
$ws1.Range("F13").Value = 96

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 97
$ws4.Range("F3").Value = 823
$ws4.Range("F6").Value = 136
$ws4.Range("F9").Value = 4856
$ws4.Range("F11").Value = 5176
$ws4.Range("F12").Value = 591
$ws4.Range("F13").Value = 1294
$ws4.Range("F14").Value = 96
